$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held a small "business group" voting-right table
# (3 data rows, with B2 using a custom Lucida Console / vertical-centered
# style). The commit separates that graph's data out and replaces it with
# the voting-right figures for the Mellat bank group (13 rows). Clear the
# old one-off formatting on B2 so it falls back to the default style.
$ws.Range("B2").ClearFormats()

# New (index, voting right) pairs for rows 2..14
$data = @(
    @(14,  0.6773),
    @(16,  0.769),
    @(50,  0.762),
    @(112, 0.6052),
    @(254, 0.6596),
    @(321, 0.8336),
    @(355, 0.3452),
    @(373, 0.432),
    @(436, 0.6425),
    @(535, 0.5295),
    @(573, 0.2514),
    @(719, 0.5718),
    @(726, 0.2563)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Move the selection to match the new cursor position left after entering
# the data (H20), and maximize the window like the saved view.
$ws.Range("H20").Select()
$wb.Windows.Item(1).WindowState = -4137
